# Actualización automática de tasas-transfi.xlsx
$wb = $excel.ActiveWorkbook

# --- Hoja1: update the "Conversión del día" note with new exchange rates ---
$ws1 = $wb.Worksheets.Item("Hoja1")
$texto = "Conversión del día 💰`n✅ Dólar paralelo: 68`n`nBinance`n✅ 1000 Bs = 3.34 = 12960.51 pesos`n✅ 12960.51 pesos = 3.31 = 948.01 Bs`n`nPromedio competencia`n✅ Tasa pesos: 20`n✅ Tasa Bs: 20`n✅ % Ganancia: 20%"
$ws1.Range("A1").Value = $texto

# --- tasas: update the N10/O10 (Binance) and N12/O12 (transfi) rate cells ---
$ws2 = $wb.Worksheets.Item("tasas")
$ws2.Range("N10").Value = 298.989
$ws2.Range("O10").Value = 3875.05
$ws2.Range("N12").Value = 3909.99
$ws2.Range("O12").Value = 286
